$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metricas")

# Row 18: new data item "Código para lectura/escritura del SEL"
$ws.Range("A18").Value = "Código para lectura/escritura del SEL"
$ws.Range("B18").Value = 70
$ws.Range("C18").Value = 86
$ws.Range("D18").Value = 0.013888888888888888
$ws.Range("E18").Value = 0.1875
$ws.Range("F18").Value = 0.21666666666666667
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0

# Update the selection to reflect the edit
$ws.Range("C24").Select()
